# Auto update: 2025-12-06 01:15:18
# Updates the latest market-data snapshot values on Sheet1 for ASML, TSM,
# AMD, NVDA and QCOM (rows 2-6): closing price, RSI, 5-day return,
# (AMD's) 10-day up-probability, final score and MACRO_SCORE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - ASML
$ws.Range("D2").Value = 1111.18
$ws.Range("E2").Value = 65.2
$ws.Range("F2").Value = 4.83
$ws.Range("K2").Value = 63.6
$ws.Range("N2").Value = 50.60178744571824

# Row 3 - TSM
$ws.Range("D3").Value = 296.63
$ws.Range("E3").Value = 61.2
$ws.Range("F3").Value = 1.76
$ws.Range("K3").Value = 58.4
$ws.Range("N3").Value = 50.60178744571824

# Row 4 - AMD
$ws.Range("D4").Value = 218.29
$ws.Range("E4").Value = 33.7
$ws.Range("F4").Value = 0.35
$ws.Range("J4").Value = 76
$ws.Range("K4").Value = 55.6
$ws.Range("N4").Value = 50.60178744571824

# Row 5 - NVDA
$ws.Range("D5").Value = 181.86
$ws.Range("E5").Value = 41.2
$ws.Range("F5").Value = 2.75
$ws.Range("K5").Value = 50.6
$ws.Range("N5").Value = 50.60178744571824

# Row 6 - QCOM
$ws.Range("D6").Value = 176.21
$ws.Range("E6").Value = 54.2
$ws.Range("F6").Value = 5.37
$ws.Range("K6").Value = 49.2
$ws.Range("N6").Value = 50.60178744571824
